$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove trailing column BA (last column) and trailing rows 23-24 (last rows)
$ws.Columns("BA").Delete()
$ws.Rows("24").Delete()
$ws.Rows("23").Delete()

# Step 2: clear cells that become blank (bugfixed forecaster no longer back-fills these early periods)
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G6").ClearContents()
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("O11").ClearContents()
$ws.Range("Q12").ClearContents()
$ws.Range("R13").ClearContents()
$ws.Range("S13").ClearContents()
$ws.Range("T14").ClearContents()
$ws.Range("U14").ClearContents()
$ws.Range("V14").ClearContents()
$ws.Range("V15").ClearContents()
$ws.Range("W15").ClearContents()
$ws.Range("X15").ClearContents()
$ws.Range("Y15").ClearContents()
$ws.Range("Z15").ClearContents()
$ws.Range("AA16").ClearContents()
$ws.Range("AB16").ClearContents()
$ws.Range("AC16").ClearContents()
$ws.Range("AD16").ClearContents()
$ws.Range("Y16").ClearContents()
$ws.Range("Z16").ClearContents()
$ws.Range("AC17").ClearContents()
$ws.Range("AD17").ClearContents()
$ws.Range("AE17").ClearContents()
$ws.Range("AF17").ClearContents()
$ws.Range("AG17").ClearContents()
$ws.Range("AG18").ClearContents()
$ws.Range("AH18").ClearContents()
$ws.Range("AI18").ClearContents()
$ws.Range("AJ18").ClearContents()
$ws.Range("AK18").ClearContents()
$ws.Range("AK19").ClearContents()
$ws.Range("AL19").ClearContents()
$ws.Range("AM19").ClearContents()
$ws.Range("AN19").ClearContents()
$ws.Range("AO19").ClearContents()
$ws.Range("AO20").ClearContents()
$ws.Range("AP20").ClearContents()
$ws.Range("AQ20").ClearContents()
$ws.Range("AR20").ClearContents()
$ws.Range("AS20").ClearContents()
$ws.Range("AS21").ClearContents()
$ws.Range("AT21").ClearContents()
$ws.Range("AU21").ClearContents()
$ws.Range("AV21").ClearContents()
$ws.Range("AW21").ClearContents()
$ws.Range("AW22").ClearContents()
$ws.Range("AX22").ClearContents()
$ws.Range("AY22").ClearContents()
$ws.Range("AZ22").ClearContents()

# Step 3: update cell values that changed (corrected forecast values + shifted header dates)
$ws.Range("AA1").Value = 43600
$ws.Range("AB1").Value = 43691
$ws.Range("AC1").Value = 43783
$ws.Range("AD1").Value = 43875
$ws.Range("AE1").Value = 43966
$ws.Range("AF1").Value = 44068
$ws.Range("AG1").Value = 44159
$ws.Range("AH1").Value = 44251
$ws.Range("AI1").Value = 44341
$ws.Range("AJ1").Value = 44432
$ws.Range("AK1").Value = 44525
$ws.Range("AL1").Value = 44617
$ws.Range("AM1").Value = 44706
$ws.Range("AN1").Value = 44798
$ws.Range("AO1").Value = 44890
$ws.Range("AP1").Value = 44981
$ws.Range("AQ1").Value = 45071
$ws.Range("AR1").Value = 45163
$ws.Range("AS1").Value = 45254
$ws.Range("AT1").Value = 45345
$ws.Range("AU1").Value = 45436
$ws.Range("AV1").Value = 45534
$ws.Range("AW1").Value = 45618
$ws.Range("AX1").Value = 45713
$ws.Range("AY1").Value = 45800
$ws.Range("AZ1").Value = 45891
$ws.Range("B1").Value = 39583
$ws.Range("C1").Value = 39765
$ws.Range("D1").Value = 39948
$ws.Range("E1").Value = 40130
$ws.Range("F1").Value = 40310
$ws.Range("G1").Value = 40494
$ws.Range("H1").Value = 40676
$ws.Range("I1").Value = 40862
$ws.Range("J1").Value = 41044
$ws.Range("K1").Value = 41228
$ws.Range("L1").Value = 41409
$ws.Range("M1").Value = 41592
$ws.Range("N1").Value = 41774
$ws.Range("O1").Value = 41957
$ws.Range("P1").Value = 42137
$ws.Range("Q1").Value = 42321
$ws.Range("R1").Value = 42503
$ws.Range("S1").Value = 42689
$ws.Range("T1").Value = 42867
$ws.Range("U1").Value = 43053
$ws.Range("V1").Value = 43145
$ws.Range("W1").Value = 43235
$ws.Range("X1").Value = 43326
$ws.Range("Y1").Value = 43418
$ws.Range("Z1").Value = 43510
$ws.Range("E3").Value = 0.1715429114845124
$ws.Range("F3").Value = 0.1715429114845124
$ws.Range("G3").Value = 0.1715429114845124
$ws.Range("H3").Value = 0.1715429114845124
$ws.Range("I3").Value = 0.1715429114845124
$ws.Range("J3").Value = 0.1715429114845124
$ws.Range("K7").Value = 0.5784444854042281
$ws.Range("L7").Value = 1.133560223479058
$ws.Range("K8").Value = 1.962049292219414
$ws.Range("L8").Value = 1.985690391709771
$ws.Range("M8").Value = 2.529895848567842
$ws.Range("N8").Value = 3.633318781899142
$ws.Range("L9").Value = 2.090102686531425
$ws.Range("M9").Value = 2.205381251914007
$ws.Range("N9").Value = 2.715291551682419
$ws.Range("O9").Value = 4.060884847379076
$ws.Range("P9").Value = 3.057638025163611
$ws.Range("N10").Value = 2.448864397591044
$ws.Range("O10").Value = 2.792143403677905
$ws.Range("P10").Value = 2.42782168586293
$ws.Range("Q10").Value = 2.270469368501771
$ws.Range("R10").Value = 2.319057151538662
$ws.Range("P11").Value = 2.543955481275106
$ws.Range("Q11").Value = 2.507859322024841
$ws.Range("R11").Value = 2.508920621023392
$ws.Range("S11").Value = 2.467161166346266
$ws.Range("T11").Value = 2.536029549059826
$ws.Range("R12").Value = 2.526834392238175
$ws.Range("S12").Value = 2.518575433256176
$ws.Range("T12").Value = 2.546671316138061
$ws.Range("U12").Value = 2.480855794925163
$ws.Range("V12").Value = 3.025024236774643
$ws.Range("W12").Value = 3.120740332206995
$ws.Range("X12").Value = 3.279355759764568
$ws.Range("AA13").Value = 2.891533899000343
$ws.Range("AB13").Value = 2.827707622797226
$ws.Range("T13").Value = 2.530440776250154
$ws.Range("U13").Value = 2.518755579319643
$ws.Range("V13").Value = 2.69389938681992
$ws.Range("W13").Value = 2.775533179497169
$ws.Range("X13").Value = 3.107596903291299
$ws.Range("Y13").Value = 3.221757900820066
$ws.Range("Z13").Value = 2.945303709067959
$ws.Range("AA14").Value = 2.545843589346886
$ws.Range("AB14").Value = 2.413544192054795
$ws.Range("AC14").Value = 2.631992339577627
$ws.Range("AD14").Value = 2.552688975800033
$ws.Range("AE14").Value = 2.618329006605924
$ws.Range("AF14").Value = 1.790319754067715
$ws.Range("W14").Value = 2.604201945499174
$ws.Range("X14").Value = 2.657071530429667
$ws.Range("Y14").Value = 2.678174398932609
$ws.Range("Z14").Value = 2.591074440292807
$ws.Range("AA15").Value = 2.623024301937549
$ws.Range("AB15").Value = 2.60322048149817
$ws.Range("AC15").Value = 2.640819364776803
$ws.Range("AD15").Value = 2.616345720823721
$ws.Range("AE15").Value = 2.671430903007876
$ws.Range("AF15").Value = 1.691013991470625
$ws.Range("AG15").Value = 2.153309886824961
$ws.Range("AH15").Value = 2.130407351599706
$ws.Range("AI15").Value = 2.137626121054947
$ws.Range("AJ15").Value = 2.339531676162721
$ws.Range("AE16").Value = 2.630644791314363
$ws.Range("AF16").Value = 2.558570068847144
$ws.Range("AG16").Value = 2.761341020331276
$ws.Range("AH16").Value = 2.785334366326175
$ws.Range("AI16").Value = 2.891950990452763
$ws.Range("AJ16").Value = 3.941556826710224
$ws.Range("AK16").Value = 4.667362054855917
$ws.Range("AL16").Value = 5.037171918133976
$ws.Range("AM16").Value = 4.951039758187648
$ws.Range("AN16").Value = 4.834496776263886
$ws.Range("AH17").Value = 2.764442819703916
$ws.Range("AI17").Value = 2.919819837356252
$ws.Range("AJ17").Value = 3.134394395265594
$ws.Range("AK17").Value = 3.327089769540992
$ws.Range("AL17").Value = 3.641364543513781
$ws.Range("AM17").Value = 3.481452844954491
$ws.Range("AN17").Value = 2.845322256798233
$ws.Range("AO17").Value = 3.305715257492858
$ws.Range("AP17").Value = 3.153537734543965
$ws.Range("AQ17").Value = 2.838865660558509
$ws.Range("AR17").Value = 2.798216547494237
$ws.Range("AL18").Value = 3.226452504784616
$ws.Range("AM18").Value = 3.20116940334636
$ws.Range("AN18").Value = 3.029622899744266
$ws.Range("AO18").Value = 3.019047171689593
$ws.Range("AP18").Value = 2.935215611250452
$ws.Range("AQ18").Value = 2.377254777217375
$ws.Range("AR18").Value = 2.138412043368865
$ws.Range("AS18").Value = 1.757655717321982
$ws.Range("AT18").Value = 1.831762447564067
$ws.Range("AU18").Value = 1.625773169906108
$ws.Range("AV18").Value = 1.530879676868468
$ws.Range("AP19").Value = 2.986397903652205
$ws.Range("AQ19").Value = 2.82910658530624
$ws.Range("AR19").Value = 2.714478023861111
$ws.Range("AS19").Value = 2.633539027099796
$ws.Range("AT19").Value = 2.69124964061378
$ws.Range("AU19").Value = 2.42082970885531
$ws.Range("AV19").Value = 2.01742511619909
$ws.Range("AW19").Value = 2.159361127638926
$ws.Range("AX19").Value = 2.104676416355189
$ws.Range("AY19").Value = 2.030491763452114
$ws.Range("AZ19").Value = 2.060859685319461
$ws.Range("AT20").Value = 2.711808184127418
$ws.Range("AU20").Value = 2.79751891585911
$ws.Range("AV20").Value = 2.787508609954714
$ws.Range("AW20").Value = 2.780289798993185
$ws.Range("AX20").Value = 2.754798876280251
$ws.Range("AY20").Value = 2.559374235215039
$ws.Range("AZ20").Value = 2.733459627814305
$ws.Range("AX21").Value = 2.825169002342753
$ws.Range("AY21").Value = 2.70120649680623
$ws.Range("AZ21").Value = 2.76671919604734

Write-Output "done"